$wb = $excel.ActiveWorkbook
$worksheets = $wb.Worksheets

# --- "展览" worksheet ---
$ws1 = $worksheets.Item("展览")
$ws1.Range("F3").Value = 477
$ws1.Range("F6").Value = 323
$ws1.Range("F10").Value = 1319
$ws1.Range("F11").Value = 180
$ws1.Range("G13").Value = 28
$ws1.Range("F18").Value = 1672
$ws1.Range("F21").Value = 236
$ws1.Range("F22").Value = 2485
$ws1.Range("F23").Value = 16
$ws1.Range("F30").Value = 2832
$ws1.Range("F31").Value = 1629
$ws1.Range("F34").Value = 680
$ws1.Range("F35").Value = 867
$ws1.Range("F36").Value = 1860
$ws1.Range("F37").Value = 898
$ws1.Range("F38").Value = 1870
$ws1.Range("F39").Value = 205
$ws1.Range("F41").Value = 845
$ws1.Range("F42").Value = 43
$ws1.Range("F43").Value = 877
$ws1.Range("F45").Value = 1029
$ws1.Range("F46").Value = 103
$ws1.Range("F48").Value = 224
$ws1.Range("F49").Value = 3351

# --- "演出" worksheet ---
$ws2 = $worksheets.Item("演出")
$ws2.Range("F12").Value = 804

# --- "全部类型" worksheet ---
$ws4 = $worksheets.Item("全部类型")
$ws4.Range("F2").Value = 477
$ws4.Range("F7").Value = 323
$ws4.Range("F11").Value = 1319
$ws4.Range("F12").Value = 180
$ws4.Range("G14").Value = 28
$ws4.Range("F19").Value = 1672
$ws4.Range("F22").Value = 236
$ws4.Range("F23").Value = 2485
$ws4.Range("F27").Value = 2832
$ws4.Range("F28").Value = 1629
$ws4.Range("F32").Value = 804
$ws4.Range("F34").Value = 680
$ws4.Range("F35").Value = 867
$ws4.Range("F36").Value = 1860
$ws4.Range("F38").Value = 898
$ws4.Range("F39").Value = 1870
$ws4.Range("F40").Value = 845
$ws4.Range("F41").Value = 877
$ws4.Range("F43").Value = 1029
$ws4.Range("F44").Value = 103
$ws4.Range("F47").Value = 224
$ws4.Range("F48").Value = 3351
